$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows (1253, 1254) are appended at the end; give column D the same
# date number format used throughout the rest of the column.
$ws.Range("D1253:D1254").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the columns that stay constant for every row in this sheet (only
# needed for the 2 brand-new rows; existing rows already contain them).
$ws.Cells.Item(1253, 1).Value = 7
$ws.Cells.Item(1253, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(1253, 3).Value = 'Ñuble'
$ws.Cells.Item(1253, 5).Value = 16
$ws.Cells.Item(1253, 6).Value = 'Fruta'
$ws.Cells.Item(1253, 7).Value = 100106
$ws.Cells.Item(1253, 8).Value = 'Oleaginosos'
$ws.Cells.Item(1253, 9).Value = 100106002
$ws.Cells.Item(1253, 10).Value = 'Palta'
$ws.Cells.Item(1254, 1).Value = 7
$ws.Cells.Item(1254, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(1254, 3).Value = 'Ñuble'
$ws.Cells.Item(1254, 5).Value = 16
$ws.Cells.Item(1254, 6).Value = 'Fruta'
$ws.Cells.Item(1254, 7).Value = 100106
$ws.Cells.Item(1254, 8).Value = 'Oleaginosos'
$ws.Cells.Item(1254, 9).Value = 100106002
$ws.Cells.Item(1254, 10).Value = 'Palta'

$ws.Cells.Item(1180, 4).Value = 45265
$ws.Cells.Item(1180, 11).Value = 'Hass'
$ws.Cells.Item(1180, 12).Value = 'Especial'
$ws.Cells.Item(1180, 13).Value = 150
$ws.Cells.Item(1180, 14).Value = 3800
$ws.Cells.Item(1180, 15).Value = 3800
$ws.Cells.Item(1180, 16).Value = 3800
$ws.Cells.Item(1180, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1180, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1180, 19).Value = 3800
$ws.Cells.Item(1180, 20).Value = 1

$ws.Cells.Item(1181, 4).Value = 45265
$ws.Cells.Item(1181, 11).Value = 'Hass'
$ws.Cells.Item(1181, 12).Value = 'Primera'
$ws.Cells.Item(1181, 13).Value = 150
$ws.Cells.Item(1181, 14).Value = 3300
$ws.Cells.Item(1181, 15).Value = 3300
$ws.Cells.Item(1181, 16).Value = 3300
$ws.Cells.Item(1181, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1181, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1181, 19).Value = 3300
$ws.Cells.Item(1181, 20).Value = 1

$ws.Cells.Item(1182, 4).Value = 45222
$ws.Cells.Item(1182, 11).Value = 'Hass'
$ws.Cells.Item(1182, 12).Value = 'Especial'
$ws.Cells.Item(1182, 13).Value = 100
$ws.Cells.Item(1182, 14).Value = 3200
$ws.Cells.Item(1182, 15).Value = 3200
$ws.Cells.Item(1182, 16).Value = 3200
$ws.Cells.Item(1182, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1182, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1182, 19).Value = 3200
$ws.Cells.Item(1182, 20).Value = 1

$ws.Cells.Item(1183, 4).Value = 45222
$ws.Cells.Item(1183, 11).Value = 'Hass'
$ws.Cells.Item(1183, 12).Value = 'Primera'
$ws.Cells.Item(1183, 13).Value = 100
$ws.Cells.Item(1183, 14).Value = 3000
$ws.Cells.Item(1183, 15).Value = 3000
$ws.Cells.Item(1183, 16).Value = 3000
$ws.Cells.Item(1183, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1183, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1183, 19).Value = 3000
$ws.Cells.Item(1183, 20).Value = 1

$ws.Cells.Item(1184, 4).Value = 45222
$ws.Cells.Item(1184, 11).Value = 'Hass'
$ws.Cells.Item(1184, 12).Value = 'Segunda'
$ws.Cells.Item(1184, 13).Value = 100
$ws.Cells.Item(1184, 14).Value = 2800
$ws.Cells.Item(1184, 15).Value = 2800
$ws.Cells.Item(1184, 16).Value = 2800
$ws.Cells.Item(1184, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1184, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1184, 19).Value = 2800
$ws.Cells.Item(1184, 20).Value = 1

$ws.Cells.Item(1185, 4).Value = 45222
$ws.Cells.Item(1185, 11).Value = 'Hass'
$ws.Cells.Item(1185, 12).Value = 'Tercera'
$ws.Cells.Item(1185, 13).Value = 100
$ws.Cells.Item(1185, 14).Value = 2400
$ws.Cells.Item(1185, 15).Value = 2400
$ws.Cells.Item(1185, 16).Value = 2400
$ws.Cells.Item(1185, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1185, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1185, 19).Value = 2400
$ws.Cells.Item(1185, 20).Value = 1

$ws.Cells.Item(1186, 4).Value = 45128
$ws.Cells.Item(1186, 11).Value = 'Hass'
$ws.Cells.Item(1186, 12).Value = 'Primera'
$ws.Cells.Item(1186, 13).Value = 100
$ws.Cells.Item(1186, 14).Value = 25000
$ws.Cells.Item(1186, 15).Value = 25000
$ws.Cells.Item(1186, 16).Value = 25000
$ws.Cells.Item(1186, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1186, 18).Value = 'Perú'
$ws.Cells.Item(1186, 19).Value = 2500
$ws.Cells.Item(1186, 20).Value = 10

$ws.Cells.Item(1187, 4).Value = 45128
$ws.Cells.Item(1187, 11).Value = 'Hass'
$ws.Cells.Item(1187, 12).Value = 'Segunda'
$ws.Cells.Item(1187, 13).Value = 100
$ws.Cells.Item(1187, 14).Value = 20000
$ws.Cells.Item(1187, 15).Value = 20000
$ws.Cells.Item(1187, 16).Value = 20000
$ws.Cells.Item(1187, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1187, 18).Value = 'Perú'
$ws.Cells.Item(1187, 19).Value = 2000
$ws.Cells.Item(1187, 20).Value = 10

$ws.Cells.Item(1188, 4).Value = 44363
$ws.Cells.Item(1188, 11).Value = 'Hass'
$ws.Cells.Item(1188, 12).Value = 'Primera'
$ws.Cells.Item(1188, 13).Value = 60
$ws.Cells.Item(1188, 14).Value = 32000
$ws.Cells.Item(1188, 15).Value = 32000
$ws.Cells.Item(1188, 16).Value = 32000
$ws.Cells.Item(1188, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1188, 18).Value = 'Perú'
$ws.Cells.Item(1188, 19).Value = 3200
$ws.Cells.Item(1188, 20).Value = 10

$ws.Cells.Item(1189, 4).Value = 44363
$ws.Cells.Item(1189, 11).Value = 'Hass'
$ws.Cells.Item(1189, 12).Value = 'Segunda'
$ws.Cells.Item(1189, 13).Value = 80
$ws.Cells.Item(1189, 14).Value = 30000
$ws.Cells.Item(1189, 15).Value = 31000
$ws.Cells.Item(1189, 16).Value = 30500
$ws.Cells.Item(1189, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1189, 18).Value = 'Perú'
$ws.Cells.Item(1189, 19).Value = 3050
$ws.Cells.Item(1189, 20).Value = 10

$ws.Cells.Item(1190, 4).Value = 45097
$ws.Cells.Item(1190, 11).Value = 'Hass'
$ws.Cells.Item(1190, 12).Value = 'Primera'
$ws.Cells.Item(1190, 13).Value = 150
$ws.Cells.Item(1190, 14).Value = 30000
$ws.Cells.Item(1190, 15).Value = 30000
$ws.Cells.Item(1190, 16).Value = 30000
$ws.Cells.Item(1190, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1190, 18).Value = 'Perú'
$ws.Cells.Item(1190, 19).Value = 3000
$ws.Cells.Item(1190, 20).Value = 10

$ws.Cells.Item(1191, 4).Value = 45097
$ws.Cells.Item(1191, 11).Value = 'Hass'
$ws.Cells.Item(1191, 12).Value = 'Segunda'
$ws.Cells.Item(1191, 13).Value = 120
$ws.Cells.Item(1191, 14).Value = 25000
$ws.Cells.Item(1191, 15).Value = 25000
$ws.Cells.Item(1191, 16).Value = 25000
$ws.Cells.Item(1191, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1191, 18).Value = 'Perú'
$ws.Cells.Item(1191, 19).Value = 2500
$ws.Cells.Item(1191, 20).Value = 10

$ws.Cells.Item(1192, 4).Value = 44222
$ws.Cells.Item(1192, 11).Value = 'Hass'
$ws.Cells.Item(1192, 12).Value = 'Especial'
$ws.Cells.Item(1192, 13).Value = 42
$ws.Cells.Item(1192, 14).Value = 4900
$ws.Cells.Item(1192, 15).Value = 5000
$ws.Cells.Item(1192, 16).Value = 4948
$ws.Cells.Item(1192, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1192, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1192, 19).Value = 4948
$ws.Cells.Item(1192, 20).Value = 1

$ws.Cells.Item(1193, 4).Value = 44222
$ws.Cells.Item(1193, 11).Value = 'Hass'
$ws.Cells.Item(1193, 12).Value = 'Primera'
$ws.Cells.Item(1193, 13).Value = 60
$ws.Cells.Item(1193, 14).Value = 4700
$ws.Cells.Item(1193, 15).Value = 4800
$ws.Cells.Item(1193, 16).Value = 4758
$ws.Cells.Item(1193, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1193, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1193, 19).Value = 4758
$ws.Cells.Item(1193, 20).Value = 1

$ws.Cells.Item(1194, 4).Value = 44222
$ws.Cells.Item(1194, 11).Value = 'Hass'
$ws.Cells.Item(1194, 12).Value = 'Segunda'
$ws.Cells.Item(1194, 13).Value = 48
$ws.Cells.Item(1194, 14).Value = 4500
$ws.Cells.Item(1194, 15).Value = 4600
$ws.Cells.Item(1194, 16).Value = 4552
$ws.Cells.Item(1194, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1194, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1194, 19).Value = 4552
$ws.Cells.Item(1194, 20).Value = 1

$ws.Cells.Item(1195, 4).Value = 44271
$ws.Cells.Item(1195, 11).Value = 'Hass'
$ws.Cells.Item(1195, 12).Value = 'Primera'
$ws.Cells.Item(1195, 13).Value = 60
$ws.Cells.Item(1195, 14).Value = 4900
$ws.Cells.Item(1195, 15).Value = 5000
$ws.Cells.Item(1195, 16).Value = 4950
$ws.Cells.Item(1195, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1195, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1195, 19).Value = 4950
$ws.Cells.Item(1195, 20).Value = 1

$ws.Cells.Item(1196, 4).Value = 44271
$ws.Cells.Item(1196, 11).Value = 'Hass'
$ws.Cells.Item(1196, 12).Value = 'Segunda'
$ws.Cells.Item(1196, 13).Value = 20
$ws.Cells.Item(1196, 14).Value = 4600
$ws.Cells.Item(1196, 15).Value = 4600
$ws.Cells.Item(1196, 16).Value = 4600
$ws.Cells.Item(1196, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1196, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1196, 19).Value = 4600
$ws.Cells.Item(1196, 20).Value = 1

$ws.Cells.Item(1197, 4).Value = 44880
$ws.Cells.Item(1197, 11).Value = 'Hass'
$ws.Cells.Item(1197, 12).Value = '1a nueva(o)'
$ws.Cells.Item(1197, 13).Value = 100
$ws.Cells.Item(1197, 14).Value = 2600
$ws.Cells.Item(1197, 15).Value = 2700
$ws.Cells.Item(1197, 16).Value = 2650
$ws.Cells.Item(1197, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1197, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1197, 19).Value = 2650
$ws.Cells.Item(1197, 20).Value = 1

$ws.Cells.Item(1198, 4).Value = 44880
$ws.Cells.Item(1198, 11).Value = 'Hass'
$ws.Cells.Item(1198, 12).Value = '2a nueva(o)'
$ws.Cells.Item(1198, 13).Value = 120
$ws.Cells.Item(1198, 14).Value = 2400
$ws.Cells.Item(1198, 15).Value = 2500
$ws.Cells.Item(1198, 16).Value = 2450
$ws.Cells.Item(1198, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1198, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1198, 19).Value = 2450
$ws.Cells.Item(1198, 20).Value = 1

$ws.Cells.Item(1199, 4).Value = 44399
$ws.Cells.Item(1199, 11).Value = 'Hass'
$ws.Cells.Item(1199, 12).Value = 'Primera'
$ws.Cells.Item(1199, 13).Value = 120
$ws.Cells.Item(1199, 14).Value = 24000
$ws.Cells.Item(1199, 15).Value = 25000
$ws.Cells.Item(1199, 16).Value = 24500
$ws.Cells.Item(1199, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1199, 18).Value = 'Perú'
$ws.Cells.Item(1199, 19).Value = 2450
$ws.Cells.Item(1199, 20).Value = 10

$ws.Cells.Item(1200, 4).Value = 44435
$ws.Cells.Item(1200, 11).Value = 'Hass'
$ws.Cells.Item(1200, 12).Value = 'Primera'
$ws.Cells.Item(1200, 13).Value = 2080
$ws.Cells.Item(1200, 14).Value = 20000
$ws.Cells.Item(1200, 15).Value = 25000
$ws.Cells.Item(1200, 16).Value = 22192
$ws.Cells.Item(1200, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1200, 18).Value = 'Perú'
$ws.Cells.Item(1200, 19).Value = 2219
$ws.Cells.Item(1200, 20).Value = 10

$ws.Cells.Item(1201, 4).Value = 44435
$ws.Cells.Item(1201, 11).Value = 'Hass'
$ws.Cells.Item(1201, 12).Value = 'Segunda'
$ws.Cells.Item(1201, 13).Value = 1400
$ws.Cells.Item(1201, 14).Value = 18000
$ws.Cells.Item(1201, 15).Value = 23000
$ws.Cells.Item(1201, 16).Value = 20386
$ws.Cells.Item(1201, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1201, 18).Value = 'Perú'
$ws.Cells.Item(1201, 19).Value = 2039
$ws.Cells.Item(1201, 20).Value = 10

$ws.Cells.Item(1202, 4).Value = 44509
$ws.Cells.Item(1202, 11).Value = 'Hass'
$ws.Cells.Item(1202, 12).Value = 'Primera'
$ws.Cells.Item(1202, 13).Value = 120
$ws.Cells.Item(1202, 14).Value = 2800
$ws.Cells.Item(1202, 15).Value = 2900
$ws.Cells.Item(1202, 16).Value = 2850
$ws.Cells.Item(1202, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1202, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1202, 19).Value = 2850
$ws.Cells.Item(1202, 20).Value = 1

$ws.Cells.Item(1203, 4).Value = 44509
$ws.Cells.Item(1203, 11).Value = 'Hass'
$ws.Cells.Item(1203, 12).Value = 'Segunda'
$ws.Cells.Item(1203, 13).Value = 120
$ws.Cells.Item(1203, 14).Value = 2500
$ws.Cells.Item(1203, 15).Value = 2600
$ws.Cells.Item(1203, 16).Value = 2550
$ws.Cells.Item(1203, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1203, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1203, 19).Value = 2550
$ws.Cells.Item(1203, 20).Value = 1

$ws.Cells.Item(1204, 4).Value = 45050
$ws.Cells.Item(1204, 11).Value = 'Hass'
$ws.Cells.Item(1204, 12).Value = 'Especial'
$ws.Cells.Item(1204, 13).Value = 100
$ws.Cells.Item(1204, 14).Value = 30000
$ws.Cells.Item(1204, 15).Value = 35000
$ws.Cells.Item(1204, 16).Value = 32500
$ws.Cells.Item(1204, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1204, 18).Value = 'Perú'
$ws.Cells.Item(1204, 19).Value = 3250
$ws.Cells.Item(1204, 20).Value = 10

$ws.Cells.Item(1205, 4).Value = 45050
$ws.Cells.Item(1205, 11).Value = 'Hass'
$ws.Cells.Item(1205, 12).Value = 'Primera'
$ws.Cells.Item(1205, 13).Value = 80
$ws.Cells.Item(1205, 14).Value = 28000
$ws.Cells.Item(1205, 15).Value = 28000
$ws.Cells.Item(1205, 16).Value = 28000
$ws.Cells.Item(1205, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1205, 18).Value = 'Perú'
$ws.Cells.Item(1205, 19).Value = 2800
$ws.Cells.Item(1205, 20).Value = 10

$ws.Cells.Item(1206, 4).Value = 45050
$ws.Cells.Item(1206, 11).Value = 'Hass'
$ws.Cells.Item(1206, 12).Value = 'Segunda'
$ws.Cells.Item(1206, 13).Value = 80
$ws.Cells.Item(1206, 14).Value = 25000
$ws.Cells.Item(1206, 15).Value = 25000
$ws.Cells.Item(1206, 16).Value = 25000
$ws.Cells.Item(1206, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1206, 18).Value = 'Perú'
$ws.Cells.Item(1206, 19).Value = 2500
$ws.Cells.Item(1206, 20).Value = 10

$ws.Cells.Item(1207, 4).Value = 45180
$ws.Cells.Item(1207, 11).Value = 'Hass'
$ws.Cells.Item(1207, 12).Value = '1a nueva(o)'
$ws.Cells.Item(1207, 13).Value = 100
$ws.Cells.Item(1207, 14).Value = 3500
$ws.Cells.Item(1207, 15).Value = 3500
$ws.Cells.Item(1207, 16).Value = 3500
$ws.Cells.Item(1207, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1207, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1207, 19).Value = 3500
$ws.Cells.Item(1207, 20).Value = 1

$ws.Cells.Item(1208, 4).Value = 45180
$ws.Cells.Item(1208, 11).Value = 'Hass'
$ws.Cells.Item(1208, 12).Value = '2a nueva(o)'
$ws.Cells.Item(1208, 13).Value = 100
$ws.Cells.Item(1208, 14).Value = 3000
$ws.Cells.Item(1208, 15).Value = 3000
$ws.Cells.Item(1208, 16).Value = 3000
$ws.Cells.Item(1208, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1208, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1208, 19).Value = 3000
$ws.Cells.Item(1208, 20).Value = 1

$ws.Cells.Item(1209, 4).Value = 45180
$ws.Cells.Item(1209, 11).Value = 'Hass'
$ws.Cells.Item(1209, 12).Value = 'Especial'
$ws.Cells.Item(1209, 13).Value = 100
$ws.Cells.Item(1209, 14).Value = 30000
$ws.Cells.Item(1209, 15).Value = 30000
$ws.Cells.Item(1209, 16).Value = 30000
$ws.Cells.Item(1209, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1209, 18).Value = 'Perú'
$ws.Cells.Item(1209, 19).Value = 3000
$ws.Cells.Item(1209, 20).Value = 10

$ws.Cells.Item(1210, 4).Value = 45180
$ws.Cells.Item(1210, 11).Value = 'Hass'
$ws.Cells.Item(1210, 12).Value = 'Primera'
$ws.Cells.Item(1210, 13).Value = 100
$ws.Cells.Item(1210, 14).Value = 25000
$ws.Cells.Item(1210, 15).Value = 25000
$ws.Cells.Item(1210, 16).Value = 25000
$ws.Cells.Item(1210, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1210, 18).Value = 'Perú'
$ws.Cells.Item(1210, 19).Value = 2500
$ws.Cells.Item(1210, 20).Value = 10

$ws.Cells.Item(1211, 4).Value = 45180
$ws.Cells.Item(1211, 11).Value = 'Hass'
$ws.Cells.Item(1211, 12).Value = 'Segunda'
$ws.Cells.Item(1211, 13).Value = 80
$ws.Cells.Item(1211, 14).Value = 22000
$ws.Cells.Item(1211, 15).Value = 22000
$ws.Cells.Item(1211, 16).Value = 22000
$ws.Cells.Item(1211, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1211, 18).Value = 'Perú'
$ws.Cells.Item(1211, 19).Value = 2200
$ws.Cells.Item(1211, 20).Value = 10

$ws.Cells.Item(1212, 4).Value = 44824
$ws.Cells.Item(1212, 11).Value = 'Hass'
$ws.Cells.Item(1212, 12).Value = 'Primera'
$ws.Cells.Item(1212, 13).Value = 120
$ws.Cells.Item(1212, 14).Value = 22000
$ws.Cells.Item(1212, 15).Value = 23000
$ws.Cells.Item(1212, 16).Value = 22500
$ws.Cells.Item(1212, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1212, 18).Value = 'Perú'
$ws.Cells.Item(1212, 19).Value = 2250
$ws.Cells.Item(1212, 20).Value = 10

$ws.Cells.Item(1213, 4).Value = 44824
$ws.Cells.Item(1213, 11).Value = 'Hass'
$ws.Cells.Item(1213, 12).Value = 'Segunda'
$ws.Cells.Item(1213, 13).Value = 120
$ws.Cells.Item(1213, 14).Value = 20000
$ws.Cells.Item(1213, 15).Value = 21000
$ws.Cells.Item(1213, 16).Value = 20500
$ws.Cells.Item(1213, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1213, 18).Value = 'Perú'
$ws.Cells.Item(1213, 19).Value = 2050
$ws.Cells.Item(1213, 20).Value = 10

$ws.Cells.Item(1214, 4).Value = 44539
$ws.Cells.Item(1214, 11).Value = 'Hass'
$ws.Cells.Item(1214, 12).Value = 'Primera'
$ws.Cells.Item(1214, 13).Value = 500
$ws.Cells.Item(1214, 14).Value = 2700
$ws.Cells.Item(1214, 15).Value = 2800
$ws.Cells.Item(1214, 16).Value = 2750
$ws.Cells.Item(1214, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1214, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1214, 19).Value = 2750
$ws.Cells.Item(1214, 20).Value = 1

$ws.Cells.Item(1215, 4).Value = 44539
$ws.Cells.Item(1215, 11).Value = 'Hass'
$ws.Cells.Item(1215, 12).Value = 'Segunda'
$ws.Cells.Item(1215, 13).Value = 300
$ws.Cells.Item(1215, 14).Value = 2500
$ws.Cells.Item(1215, 15).Value = 2600
$ws.Cells.Item(1215, 16).Value = 2550
$ws.Cells.Item(1215, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1215, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1215, 19).Value = 2550
$ws.Cells.Item(1215, 20).Value = 1

$ws.Cells.Item(1216, 4).Value = 44894
$ws.Cells.Item(1216, 11).Value = 'Hass'
$ws.Cells.Item(1216, 12).Value = 'Primera'
$ws.Cells.Item(1216, 13).Value = 120
$ws.Cells.Item(1216, 14).Value = 2700
$ws.Cells.Item(1216, 15).Value = 2800
$ws.Cells.Item(1216, 16).Value = 2750
$ws.Cells.Item(1216, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1216, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1216, 19).Value = 2750
$ws.Cells.Item(1216, 20).Value = 1

$ws.Cells.Item(1217, 4).Value = 44894
$ws.Cells.Item(1217, 11).Value = 'Hass'
$ws.Cells.Item(1217, 12).Value = 'Segunda'
$ws.Cells.Item(1217, 13).Value = 60
$ws.Cells.Item(1217, 14).Value = 2500
$ws.Cells.Item(1217, 15).Value = 2500
$ws.Cells.Item(1217, 16).Value = 2500
$ws.Cells.Item(1217, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1217, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1217, 19).Value = 2500
$ws.Cells.Item(1217, 20).Value = 1

$ws.Cells.Item(1218, 4).Value = 44917
$ws.Cells.Item(1218, 11).Value = 'Hass'
$ws.Cells.Item(1218, 12).Value = 'Primera'
$ws.Cells.Item(1218, 13).Value = 60
$ws.Cells.Item(1218, 14).Value = 3000
$ws.Cells.Item(1218, 15).Value = 3200
$ws.Cells.Item(1218, 16).Value = 3100
$ws.Cells.Item(1218, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1218, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1218, 19).Value = 310
$ws.Cells.Item(1218, 20).Value = 10

$ws.Cells.Item(1219, 4).Value = 44917
$ws.Cells.Item(1219, 11).Value = 'Hass'
$ws.Cells.Item(1219, 12).Value = 'Segunda'
$ws.Cells.Item(1219, 13).Value = 30
$ws.Cells.Item(1219, 14).Value = 3000
$ws.Cells.Item(1219, 15).Value = 3000
$ws.Cells.Item(1219, 16).Value = 3000
$ws.Cells.Item(1219, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1219, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1219, 19).Value = 300
$ws.Cells.Item(1219, 20).Value = 10

$ws.Cells.Item(1220, 4).Value = 44914
$ws.Cells.Item(1220, 11).Value = 'Mexicola'
$ws.Cells.Item(1220, 12).Value = 'Primera'
$ws.Cells.Item(1220, 13).Value = 120
$ws.Cells.Item(1220, 14).Value = 3100
$ws.Cells.Item(1220, 15).Value = 3200
$ws.Cells.Item(1220, 16).Value = 3150
$ws.Cells.Item(1220, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1220, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1220, 19).Value = 3150
$ws.Cells.Item(1220, 20).Value = 1

$ws.Cells.Item(1221, 4).Value = 45168
$ws.Cells.Item(1221, 11).Value = 'Hass'
$ws.Cells.Item(1221, 12).Value = 'Especial'
$ws.Cells.Item(1221, 13).Value = 80
$ws.Cells.Item(1221, 14).Value = 32000
$ws.Cells.Item(1221, 15).Value = 32000
$ws.Cells.Item(1221, 16).Value = 32000
$ws.Cells.Item(1221, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1221, 18).Value = 'Perú'
$ws.Cells.Item(1221, 19).Value = 3200
$ws.Cells.Item(1221, 20).Value = 10

$ws.Cells.Item(1222, 4).Value = 45168
$ws.Cells.Item(1222, 11).Value = 'Hass'
$ws.Cells.Item(1222, 12).Value = 'Primera'
$ws.Cells.Item(1222, 13).Value = 100
$ws.Cells.Item(1222, 14).Value = 27000
$ws.Cells.Item(1222, 15).Value = 27000
$ws.Cells.Item(1222, 16).Value = 27000
$ws.Cells.Item(1222, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1222, 18).Value = 'Perú'
$ws.Cells.Item(1222, 19).Value = 2700
$ws.Cells.Item(1222, 20).Value = 10

$ws.Cells.Item(1223, 4).Value = 45168
$ws.Cells.Item(1223, 11).Value = 'Hass'
$ws.Cells.Item(1223, 12).Value = 'Segunda'
$ws.Cells.Item(1223, 13).Value = 100
$ws.Cells.Item(1223, 14).Value = 25000
$ws.Cells.Item(1223, 15).Value = 25000
$ws.Cells.Item(1223, 16).Value = 25000
$ws.Cells.Item(1223, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1223, 18).Value = 'Perú'
$ws.Cells.Item(1223, 19).Value = 2500
$ws.Cells.Item(1223, 20).Value = 10

$ws.Cells.Item(1224, 4).Value = 44574
$ws.Cells.Item(1224, 11).Value = 'Hass'
$ws.Cells.Item(1224, 12).Value = 'Primera'
$ws.Cells.Item(1224, 13).Value = 60
$ws.Cells.Item(1224, 14).Value = 2700
$ws.Cells.Item(1224, 15).Value = 2800
$ws.Cells.Item(1224, 16).Value = 2750
$ws.Cells.Item(1224, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1224, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1224, 19).Value = 2750
$ws.Cells.Item(1224, 20).Value = 1

$ws.Cells.Item(1225, 4).Value = 44574
$ws.Cells.Item(1225, 11).Value = 'Hass'
$ws.Cells.Item(1225, 12).Value = 'Segunda'
$ws.Cells.Item(1225, 13).Value = 120
$ws.Cells.Item(1225, 14).Value = 2500
$ws.Cells.Item(1225, 15).Value = 2600
$ws.Cells.Item(1225, 16).Value = 2550
$ws.Cells.Item(1225, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1225, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1225, 19).Value = 2550
$ws.Cells.Item(1225, 20).Value = 1

$ws.Cells.Item(1226, 4).Value = 44214
$ws.Cells.Item(1226, 11).Value = 'Hass'
$ws.Cells.Item(1226, 12).Value = 'Especial'
$ws.Cells.Item(1226, 13).Value = 43
$ws.Cells.Item(1226, 14).Value = 4900
$ws.Cells.Item(1226, 15).Value = 5000
$ws.Cells.Item(1226, 16).Value = 4958
$ws.Cells.Item(1226, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1226, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1226, 19).Value = 4958
$ws.Cells.Item(1226, 20).Value = 1

$ws.Cells.Item(1227, 4).Value = 44214
$ws.Cells.Item(1227, 11).Value = 'Hass'
$ws.Cells.Item(1227, 12).Value = 'Primera'
$ws.Cells.Item(1227, 13).Value = 50
$ws.Cells.Item(1227, 14).Value = 4700
$ws.Cells.Item(1227, 15).Value = 4800
$ws.Cells.Item(1227, 16).Value = 4760
$ws.Cells.Item(1227, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1227, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1227, 19).Value = 4760
$ws.Cells.Item(1227, 20).Value = 1

$ws.Cells.Item(1228, 4).Value = 44214
$ws.Cells.Item(1228, 11).Value = 'Hass'
$ws.Cells.Item(1228, 12).Value = 'Segunda'
$ws.Cells.Item(1228, 13).Value = 52
$ws.Cells.Item(1228, 14).Value = 4500
$ws.Cells.Item(1228, 15).Value = 4600
$ws.Cells.Item(1228, 16).Value = 4558
$ws.Cells.Item(1228, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1228, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1228, 19).Value = 4558
$ws.Cells.Item(1228, 20).Value = 1

$ws.Cells.Item(1229, 4).Value = 44895
$ws.Cells.Item(1229, 11).Value = 'Hass'
$ws.Cells.Item(1229, 12).Value = 'Primera'
$ws.Cells.Item(1229, 13).Value = 100
$ws.Cells.Item(1229, 14).Value = 2800
$ws.Cells.Item(1229, 15).Value = 2900
$ws.Cells.Item(1229, 16).Value = 2850
$ws.Cells.Item(1229, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1229, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1229, 19).Value = 2850
$ws.Cells.Item(1229, 20).Value = 1

$ws.Cells.Item(1230, 4).Value = 44895
$ws.Cells.Item(1230, 11).Value = 'Hass'
$ws.Cells.Item(1230, 12).Value = 'Segunda'
$ws.Cells.Item(1230, 13).Value = 50
$ws.Cells.Item(1230, 14).Value = 2600
$ws.Cells.Item(1230, 15).Value = 2600
$ws.Cells.Item(1230, 16).Value = 2600
$ws.Cells.Item(1230, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1230, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1230, 19).Value = 2600
$ws.Cells.Item(1230, 20).Value = 1

$ws.Cells.Item(1231, 4).Value = 44895
$ws.Cells.Item(1231, 11).Value = 'Hass'
$ws.Cells.Item(1231, 12).Value = 'Primera'
$ws.Cells.Item(1231, 13).Value = 120
$ws.Cells.Item(1231, 14).Value = 2700
$ws.Cells.Item(1231, 15).Value = 2800
$ws.Cells.Item(1231, 16).Value = 2750
$ws.Cells.Item(1231, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1231, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1231, 19).Value = 2750
$ws.Cells.Item(1231, 20).Value = 1

$ws.Cells.Item(1232, 4).Value = 44895
$ws.Cells.Item(1232, 11).Value = 'Hass'
$ws.Cells.Item(1232, 12).Value = 'Segunda'
$ws.Cells.Item(1232, 13).Value = 60
$ws.Cells.Item(1232, 14).Value = 2500
$ws.Cells.Item(1232, 15).Value = 2500
$ws.Cells.Item(1232, 16).Value = 2500
$ws.Cells.Item(1232, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1232, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1232, 19).Value = 2500
$ws.Cells.Item(1232, 20).Value = 1

$ws.Cells.Item(1233, 4).Value = 45135
$ws.Cells.Item(1233, 11).Value = 'Hass'
$ws.Cells.Item(1233, 12).Value = 'Especial'
$ws.Cells.Item(1233, 13).Value = 60
$ws.Cells.Item(1233, 14).Value = 28000
$ws.Cells.Item(1233, 15).Value = 28000
$ws.Cells.Item(1233, 16).Value = 28000
$ws.Cells.Item(1233, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1233, 18).Value = 'Perú'
$ws.Cells.Item(1233, 19).Value = 2800
$ws.Cells.Item(1233, 20).Value = 10

$ws.Cells.Item(1234, 4).Value = 45135
$ws.Cells.Item(1234, 11).Value = 'Hass'
$ws.Cells.Item(1234, 12).Value = 'Primera'
$ws.Cells.Item(1234, 13).Value = 60
$ws.Cells.Item(1234, 14).Value = 25000
$ws.Cells.Item(1234, 15).Value = 25000
$ws.Cells.Item(1234, 16).Value = 25000
$ws.Cells.Item(1234, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1234, 18).Value = 'Perú'
$ws.Cells.Item(1234, 19).Value = 2500
$ws.Cells.Item(1234, 20).Value = 10

$ws.Cells.Item(1235, 4).Value = 45135
$ws.Cells.Item(1235, 11).Value = 'Hass'
$ws.Cells.Item(1235, 12).Value = 'Segunda'
$ws.Cells.Item(1235, 13).Value = 60
$ws.Cells.Item(1235, 14).Value = 22000
$ws.Cells.Item(1235, 15).Value = 22000
$ws.Cells.Item(1235, 16).Value = 22000
$ws.Cells.Item(1235, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1235, 18).Value = 'Perú'
$ws.Cells.Item(1235, 19).Value = 2200
$ws.Cells.Item(1235, 20).Value = 10

$ws.Cells.Item(1236, 4).Value = 45135
$ws.Cells.Item(1236, 11).Value = 'Hass'
$ws.Cells.Item(1236, 12).Value = 'Tercera'
$ws.Cells.Item(1236, 13).Value = 120
$ws.Cells.Item(1236, 14).Value = 18000
$ws.Cells.Item(1236, 15).Value = 18000
$ws.Cells.Item(1236, 16).Value = 18000
$ws.Cells.Item(1236, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1236, 18).Value = 'Perú'
$ws.Cells.Item(1236, 19).Value = 1800
$ws.Cells.Item(1236, 20).Value = 10

$ws.Cells.Item(1237, 4).Value = 44567
$ws.Cells.Item(1237, 11).Value = 'Hass'
$ws.Cells.Item(1237, 12).Value = 'Primera'
$ws.Cells.Item(1237, 13).Value = 120
$ws.Cells.Item(1237, 14).Value = 2600
$ws.Cells.Item(1237, 15).Value = 2700
$ws.Cells.Item(1237, 16).Value = 2650
$ws.Cells.Item(1237, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1237, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1237, 19).Value = 2650
$ws.Cells.Item(1237, 20).Value = 1

$ws.Cells.Item(1238, 4).Value = 44567
$ws.Cells.Item(1238, 11).Value = 'Hass'
$ws.Cells.Item(1238, 12).Value = 'Segunda'
$ws.Cells.Item(1238, 13).Value = 120
$ws.Cells.Item(1238, 14).Value = 2400
$ws.Cells.Item(1238, 15).Value = 2500
$ws.Cells.Item(1238, 16).Value = 2450
$ws.Cells.Item(1238, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1238, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1238, 19).Value = 2450
$ws.Cells.Item(1238, 20).Value = 1

$ws.Cells.Item(1239, 4).Value = 44517
$ws.Cells.Item(1239, 11).Value = 'Hass'
$ws.Cells.Item(1239, 12).Value = 'Primera'
$ws.Cells.Item(1239, 13).Value = 160
$ws.Cells.Item(1239, 14).Value = 2600
$ws.Cells.Item(1239, 15).Value = 2700
$ws.Cells.Item(1239, 16).Value = 2650
$ws.Cells.Item(1239, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1239, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1239, 19).Value = 2650
$ws.Cells.Item(1239, 20).Value = 1

$ws.Cells.Item(1240, 4).Value = 44517
$ws.Cells.Item(1240, 11).Value = 'Hass'
$ws.Cells.Item(1240, 12).Value = 'Segunda'
$ws.Cells.Item(1240, 13).Value = 120
$ws.Cells.Item(1240, 14).Value = 2400
$ws.Cells.Item(1240, 15).Value = 2500
$ws.Cells.Item(1240, 16).Value = 2450
$ws.Cells.Item(1240, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1240, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1240, 19).Value = 2450
$ws.Cells.Item(1240, 20).Value = 1

$ws.Cells.Item(1241, 4).Value = 45121
$ws.Cells.Item(1241, 11).Value = 'Hass'
$ws.Cells.Item(1241, 12).Value = 'Primera'
$ws.Cells.Item(1241, 13).Value = 80
$ws.Cells.Item(1241, 14).Value = 25000
$ws.Cells.Item(1241, 15).Value = 25000
$ws.Cells.Item(1241, 16).Value = 25000
$ws.Cells.Item(1241, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1241, 18).Value = 'Perú'
$ws.Cells.Item(1241, 19).Value = 2500
$ws.Cells.Item(1241, 20).Value = 10

$ws.Cells.Item(1242, 4).Value = 45121
$ws.Cells.Item(1242, 11).Value = 'Hass'
$ws.Cells.Item(1242, 12).Value = 'Primera'
$ws.Cells.Item(1242, 13).Value = 60
$ws.Cells.Item(1242, 14).Value = 5000
$ws.Cells.Item(1242, 15).Value = 5000
$ws.Cells.Item(1242, 16).Value = 5000
$ws.Cells.Item(1242, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1242, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1242, 19).Value = 5000
$ws.Cells.Item(1242, 20).Value = 1

$ws.Cells.Item(1243, 4).Value = 45121
$ws.Cells.Item(1243, 11).Value = 'Hass'
$ws.Cells.Item(1243, 12).Value = 'Segunda'
$ws.Cells.Item(1243, 13).Value = 100
$ws.Cells.Item(1243, 14).Value = 20000
$ws.Cells.Item(1243, 15).Value = 20000
$ws.Cells.Item(1243, 16).Value = 20000
$ws.Cells.Item(1243, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1243, 18).Value = 'Perú'
$ws.Cells.Item(1243, 19).Value = 2000
$ws.Cells.Item(1243, 20).Value = 10

$ws.Cells.Item(1244, 4).Value = 45121
$ws.Cells.Item(1244, 11).Value = 'Hass'
$ws.Cells.Item(1244, 12).Value = 'Segunda'
$ws.Cells.Item(1244, 13).Value = 60
$ws.Cells.Item(1244, 14).Value = 4500
$ws.Cells.Item(1244, 15).Value = 4500
$ws.Cells.Item(1244, 16).Value = 4500
$ws.Cells.Item(1244, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1244, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1244, 19).Value = 4500
$ws.Cells.Item(1244, 20).Value = 1

$ws.Cells.Item(1245, 4).Value = 45121
$ws.Cells.Item(1245, 11).Value = 'Hass'
$ws.Cells.Item(1245, 12).Value = 'Tercera'
$ws.Cells.Item(1245, 13).Value = 80
$ws.Cells.Item(1245, 14).Value = 17000
$ws.Cells.Item(1245, 15).Value = 17000
$ws.Cells.Item(1245, 16).Value = 17000
$ws.Cells.Item(1245, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1245, 18).Value = 'Perú'
$ws.Cells.Item(1245, 19).Value = 1700
$ws.Cells.Item(1245, 20).Value = 10

$ws.Cells.Item(1246, 4).Value = 44818
$ws.Cells.Item(1246, 11).Value = 'Hass'
$ws.Cells.Item(1246, 12).Value = 'Primera'
$ws.Cells.Item(1246, 13).Value = 120
$ws.Cells.Item(1246, 14).Value = 2400
$ws.Cells.Item(1246, 15).Value = 2500
$ws.Cells.Item(1246, 16).Value = 2450
$ws.Cells.Item(1246, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1246, 18).Value = 'Perú'
$ws.Cells.Item(1246, 19).Value = 245
$ws.Cells.Item(1246, 20).Value = 10

$ws.Cells.Item(1247, 4).Value = 44293
$ws.Cells.Item(1247, 11).Value = 'Hass'
$ws.Cells.Item(1247, 12).Value = 'Primera'
$ws.Cells.Item(1247, 13).Value = 60
$ws.Cells.Item(1247, 14).Value = 5600
$ws.Cells.Item(1247, 15).Value = 5700
$ws.Cells.Item(1247, 16).Value = 5650
$ws.Cells.Item(1247, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1247, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1247, 19).Value = 5650
$ws.Cells.Item(1247, 20).Value = 1

$ws.Cells.Item(1248, 4).Value = 44293
$ws.Cells.Item(1248, 11).Value = 'Hass'
$ws.Cells.Item(1248, 12).Value = 'Segunda'
$ws.Cells.Item(1248, 13).Value = 60
$ws.Cells.Item(1248, 14).Value = 5100
$ws.Cells.Item(1248, 15).Value = 5200
$ws.Cells.Item(1248, 16).Value = 5150
$ws.Cells.Item(1248, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(1248, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1248, 19).Value = 5150
$ws.Cells.Item(1248, 20).Value = 1

$ws.Cells.Item(1249, 4).Value = 45100
$ws.Cells.Item(1249, 11).Value = 'Hass'
$ws.Cells.Item(1249, 12).Value = 'Especial'
$ws.Cells.Item(1249, 13).Value = 100
$ws.Cells.Item(1249, 14).Value = 35000
$ws.Cells.Item(1249, 15).Value = 35000
$ws.Cells.Item(1249, 16).Value = 35000
$ws.Cells.Item(1249, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1249, 18).Value = 'Perú'
$ws.Cells.Item(1249, 19).Value = 3500
$ws.Cells.Item(1249, 20).Value = 10

$ws.Cells.Item(1250, 4).Value = 45100
$ws.Cells.Item(1250, 11).Value = 'Hass'
$ws.Cells.Item(1250, 12).Value = 'Especial'
$ws.Cells.Item(1250, 13).Value = 120
$ws.Cells.Item(1250, 14).Value = 5000
$ws.Cells.Item(1250, 15).Value = 5000
$ws.Cells.Item(1250, 16).Value = 5000
$ws.Cells.Item(1250, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1250, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1250, 19).Value = 5000
$ws.Cells.Item(1250, 20).Value = 1

$ws.Cells.Item(1251, 4).Value = 45100
$ws.Cells.Item(1251, 11).Value = 'Hass'
$ws.Cells.Item(1251, 12).Value = 'Primera'
$ws.Cells.Item(1251, 13).Value = 80
$ws.Cells.Item(1251, 14).Value = 30000
$ws.Cells.Item(1251, 15).Value = 30000
$ws.Cells.Item(1251, 16).Value = 30000
$ws.Cells.Item(1251, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1251, 18).Value = 'Perú'
$ws.Cells.Item(1251, 19).Value = 3000
$ws.Cells.Item(1251, 20).Value = 10

$ws.Cells.Item(1252, 4).Value = 45100
$ws.Cells.Item(1252, 11).Value = 'Hass'
$ws.Cells.Item(1252, 12).Value = 'Primera'
$ws.Cells.Item(1252, 13).Value = 100
$ws.Cells.Item(1252, 14).Value = 4500
$ws.Cells.Item(1252, 15).Value = 4500
$ws.Cells.Item(1252, 16).Value = 4500
$ws.Cells.Item(1252, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1252, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1252, 19).Value = 4500
$ws.Cells.Item(1252, 20).Value = 1

$ws.Cells.Item(1253, 4).Value = 45100
$ws.Cells.Item(1253, 11).Value = 'Hass'
$ws.Cells.Item(1253, 12).Value = 'Segunda'
$ws.Cells.Item(1253, 13).Value = 100
$ws.Cells.Item(1253, 14).Value = 25000
$ws.Cells.Item(1253, 15).Value = 25000
$ws.Cells.Item(1253, 16).Value = 25000
$ws.Cells.Item(1253, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(1253, 18).Value = 'Perú'
$ws.Cells.Item(1253, 19).Value = 2500
$ws.Cells.Item(1253, 20).Value = 10

$ws.Cells.Item(1254, 4).Value = 45100
$ws.Cells.Item(1254, 11).Value = 'Hass'
$ws.Cells.Item(1254, 12).Value = 'Segunda'
$ws.Cells.Item(1254, 13).Value = 80
$ws.Cells.Item(1254, 14).Value = 4000
$ws.Cells.Item(1254, 15).Value = 4000
$ws.Cells.Item(1254, 16).Value = 4000
$ws.Cells.Item(1254, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1254, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1254, 19).Value = 4000
$ws.Cells.Item(1254, 20).Value = 1
